$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 13 - the new logbook entry
$ws.Range("C13").Value = "2/23/2021"
$ws.Range("D13").Value = "Finished border generation and movement. Began finalised sprites."
$ws.Range("E13").Value = "Border generation and movement now working as intended. There is a bug with the borders not displaying, but that will be fixed later."

# The description wraps to two lines, so the row grows taller (matches other two-line rows)
$ws.Rows.Item(13).RowHeight = 29

# Update the active selection to the description cell of the new row
$ws.Range("E13:G13").Select()
